# Adds the 2024/11/13 column (BN) of data to the "合成確率" sheet.
# Mirrors the existing per-row conditional-format colors (style 1 = no fill,
# style 2 = yellow fill / style 3 = light-blue fill) by copying formats from
# a representative cell of each existing style rather than re-deriving the rule.

$xlPasteFormats = -4122
$xlPasteValues  = -4163

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column BN needs the same width (12) as all the other date columns.
$ws.Range("BN1").ColumnWidth = 11.17

# ---- Header: BN1 = "2024/11/13" ----
# Use the same look as the other header cells (style 1 / A1) and write the
# date as literal text (via a text formula, then flattened to a value) so it
# is not auto-converted into a date serial number.
$ws.Range("A1").Copy()
$ws.Range("BN1").PasteSpecial($xlPasteFormats)
$ws.Range("BN1").Formula = '="2024/11/13"'
$ws.Range("BN1").Copy()
$ws.Range("BN1").PasteSpecial($xlPasteValues)

# ---- Data rows 2-53 ----
$ws.Range("N2").Copy()
$ws.Range("BN2").PasteSpecial($xlPasteFormats)
$ws.Range("BN2").Value = 136.7
$ws.Range("A1").Copy()
$ws.Range("BN3").PasteSpecial($xlPasteFormats)
$ws.Range("BN3").Value = 171.5
$ws.Range("A1").Copy()
$ws.Range("BN4").PasteSpecial($xlPasteFormats)
$ws.Range("BN4").Value = 230.1
$ws.Range("A1").Copy()
$ws.Range("BN5").PasteSpecial($xlPasteFormats)
$ws.Range("BN5").Value = 166.8
$ws.Range("A1").Copy()
$ws.Range("BN6").PasteSpecial($xlPasteFormats)
$ws.Range("BN6").Value = 338
$ws.Range("N2").Copy()
$ws.Range("BN7").PasteSpecial($xlPasteFormats)
$ws.Range("BN7").Value = 129
$ws.Range("N2").Copy()
$ws.Range("BN8").PasteSpecial($xlPasteFormats)
$ws.Range("BN8").Value = 125.9
$ws.Range("D2").Copy()
$ws.Range("BN9").PasteSpecial($xlPasteFormats)
$ws.Range("BN9").Value = 121.2
$ws.Range("N2").Copy()
$ws.Range("BN10").PasteSpecial($xlPasteFormats)
$ws.Range("BN10").Value = 131.3
$ws.Range("N2").Copy()
$ws.Range("BN11").PasteSpecial($xlPasteFormats)
$ws.Range("BN11").Value = 130.9
$ws.Range("A1").Copy()
$ws.Range("BN12").PasteSpecial($xlPasteFormats)
$ws.Range("BN12").Value = 152.5
$ws.Range("A1").Copy()
$ws.Range("BN13").PasteSpecial($xlPasteFormats)
$ws.Range("BN13").Value = 183.9
$ws.Range("N2").Copy()
$ws.Range("BN14").PasteSpecial($xlPasteFormats)
$ws.Range("BN14").Value = 132.9
$ws.Range("A1").Copy()
$ws.Range("BN15").PasteSpecial($xlPasteFormats)
$ws.Range("BN15").Value = 195.7
$ws.Range("A1").Copy()
$ws.Range("BN16").PasteSpecial($xlPasteFormats)
$ws.Range("BN16").Value = 150.6
$ws.Range("D2").Copy()
$ws.Range("BN17").PasteSpecial($xlPasteFormats)
$ws.Range("BN17").Value = 121
$ws.Range("D2").Copy()
$ws.Range("BN18").PasteSpecial($xlPasteFormats)
$ws.Range("BN18").Value = 107.3
$ws.Range("A1").Copy()
$ws.Range("BN19").PasteSpecial($xlPasteFormats)
$ws.Range("BN19").Value = 209.8
$ws.Range("A1").Copy()
$ws.Range("BN20").PasteSpecial($xlPasteFormats)
$ws.Range("BN20").Value = 196.3
$ws.Range("A1").Copy()
$ws.Range("BN21").PasteSpecial($xlPasteFormats)
$ws.Range("BN21").Value = 178.1
$ws.Range("N2").Copy()
$ws.Range("BN22").PasteSpecial($xlPasteFormats)
$ws.Range("BN22").Value = 131.9
$ws.Range("A1").Copy()
$ws.Range("BN23").PasteSpecial($xlPasteFormats)
$ws.Range("BN23").Value = 225.3
$ws.Range("A1").Copy()
$ws.Range("BN24").PasteSpecial($xlPasteFormats)
$ws.Range("BN24").Value = 147.9
$ws.Range("N2").Copy()
$ws.Range("BN25").PasteSpecial($xlPasteFormats)
$ws.Range("BN25").Value = 130.7
$ws.Range("N2").Copy()
$ws.Range("BN26").PasteSpecial($xlPasteFormats)
$ws.Range("BN26").Value = 134.4
$ws.Range("A1").Copy()
$ws.Range("BN27").PasteSpecial($xlPasteFormats)
$ws.Range("BN27").Value = 250.5
$ws.Range("A1").Copy()
$ws.Range("BN28").PasteSpecial($xlPasteFormats)
$ws.Range("BN28").Value = 172
$ws.Range("N2").Copy()
$ws.Range("BN29").PasteSpecial($xlPasteFormats)
$ws.Range("BN29").Value = 131.2
$ws.Range("N2").Copy()
$ws.Range("BN30").PasteSpecial($xlPasteFormats)
$ws.Range("BN30").Value = 130.2
$ws.Range("A1").Copy()
$ws.Range("BN31").PasteSpecial($xlPasteFormats)
$ws.Range("BN31").Value = 145.3
$ws.Range("A1").Copy()
$ws.Range("BN32").PasteSpecial($xlPasteFormats)
$ws.Range("BN32").Value = 142.9
$ws.Range("A1").Copy()
$ws.Range("BN33").PasteSpecial($xlPasteFormats)
$ws.Range("BN33").Value = 160.5
$ws.Range("A1").Copy()
$ws.Range("BN34").PasteSpecial($xlPasteFormats)
$ws.Range("BN34").Value = 281.3
$ws.Range("D2").Copy()
$ws.Range("BN35").PasteSpecial($xlPasteFormats)
$ws.Range("BN35").Value = 124.7
$ws.Range("A1").Copy()
$ws.Range("BN36").PasteSpecial($xlPasteFormats)
$ws.Range("BN36").Value = 196.4
$ws.Range("A1").Copy()
$ws.Range("BN37").PasteSpecial($xlPasteFormats)
$ws.Range("BN37").Value = 194.1
$ws.Range("D2").Copy()
$ws.Range("BN38").PasteSpecial($xlPasteFormats)
$ws.Range("BN38").Value = 120.6
$ws.Range("N2").Copy()
$ws.Range("BN39").PasteSpecial($xlPasteFormats)
$ws.Range("BN39").Value = 125.4
$ws.Range("D2").Copy()
$ws.Range("BN40").PasteSpecial($xlPasteFormats)
$ws.Range("BN40").Value = 119.7
$ws.Range("A1").Copy()
$ws.Range("BN41").PasteSpecial($xlPasteFormats)
$ws.Range("BN41").Value = 152.3
$ws.Range("A1").Copy()
$ws.Range("BN42").PasteSpecial($xlPasteFormats)
$ws.Range("BN42").Value = 146.6
$ws.Range("A1").Copy()
$ws.Range("BN43").PasteSpecial($xlPasteFormats)
$ws.Range("BN43").Value = 245
$ws.Range("D2").Copy()
$ws.Range("BN44").PasteSpecial($xlPasteFormats)
$ws.Range("BN44").Value = 122.8
$ws.Range("A1").Copy()
$ws.Range("BN45").PasteSpecial($xlPasteFormats)
$ws.Range("BN45").Value = 147.6
$ws.Range("A1").Copy()
$ws.Range("BN46").PasteSpecial($xlPasteFormats)
$ws.Range("BN46").Value = 141.9
$ws.Range("A1").Copy()
$ws.Range("BN47").PasteSpecial($xlPasteFormats)
$ws.Range("BN47").Value = 266.5
$ws.Range("A1").Copy()
$ws.Range("BN48").PasteSpecial($xlPasteFormats)
$ws.Range("BN48").Value = 187.6
$ws.Range("A1").Copy()
$ws.Range("BN49").PasteSpecial($xlPasteFormats)
$ws.Range("BN49").Value = 265.8
$ws.Range("A1").Copy()
$ws.Range("BN50").PasteSpecial($xlPasteFormats)
$ws.Range("BN50").Value = 161.2
$ws.Range("N2").Copy()
$ws.Range("BN51").PasteSpecial($xlPasteFormats)
$ws.Range("BN51").Value = 131.1
$ws.Range("A1").Copy()
$ws.Range("BN52").PasteSpecial($xlPasteFormats)
$ws.Range("BN52").Value = 174
$ws.Range("N2").Copy()
$ws.Range("BN53").PasteSpecial($xlPasteFormats)
$ws.Range("BN53").Value = 133.8

$excel.CutCopyMode = 0
